$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.604.04'
$ws.Range('E2').Value = '  +5.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.645.75'
$ws.Range('E3').Value = '  +11.31%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.40%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.41'
$ws.Range('E5').Value = '  +7.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.10'
$ws.Range('E6').Value = '  +13.44%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  +10.93%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.24%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('E9').Value = '  +21.32%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.52'
$ws.Range('E10').Value = '  +19.32%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.58'
$ws.Range('E11').Value = '  +4.22%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0863'
$ws.Range('E12').Value = '  +11.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.48'
$ws.Range('E13').Value = '  +22.14%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.049.79'
$ws.Range('E14').Value = '  +11.28%  '

$ws.Range('E15').Value = '  +3.98%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.659.35'
$ws.Range('E16').Value = '  +11.59%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.953'
$ws.Range('E17').Value = '  +16.08%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '15.46'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '47.768.27'
$ws.Range('E19').Value = '  +6.06%  '

$ws.Range('E20').Value = '  +12.54%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.47'
$ws.Range('E21').Value = '  +8.55%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.83'
$ws.Range('E22').Value = '  +12.52%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.17'
$ws.Range('E23').Value = '  +10.71%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '274.16'
$ws.Range('E24').Value = '  +15.36%  '

$ws.Range('E25').Value = '  +14.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.26'
$ws.Range('E26').Value = '  +20.53%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '30.80'
$ws.Range('E27').Value = '  +47.86%  '

$ws.Range('E28').Value = '  +0.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.09'
$ws.Range('E29').Value = '  +1.71%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.82'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '40.45'
$ws.Range('E31').Value = '  +8.53%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.31'
$ws.Range('E32').Value = '  +4.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.25'
$ws.Range('E33').Value = '  +16.26%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.75'
$ws.Range('E34').Value = '  -1.12%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0863'
$ws.Range('E35').Value = '  +14.11%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').Value = '  +17.59%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.89'
$ws.Range('E37').Value = '  +7.31%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.76'
$ws.Range('E38').Value = '  +3.68%  '

$ws.Range('E39').Value = '  +11.03%  '

$ws.Range('E40').Value = '  +10.38%  '

$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.68'
$ws.Range('E41').Value = '  +14.78%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.31'
$ws.Range('E42').Value = '  +56.10%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.32'
$ws.Range('E43').Value = '  +17.00%  '

$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.78'
$ws.Range('E44').Value = '  +20.16%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0336'
$ws.Range('E45').Value = '  +15.15%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.194.37'
$ws.Range('E46').Value = '  +12.19%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '96.79'
$ws.Range('E47').Value = '  +9.77%  '

$ws.Range('E48').Value = '  -0.20%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.10'
$ws.Range('E49').Value = '  +20.33%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '115.67'
$ws.Range('E50').Value = '  +16.97%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.83'
$ws.Range('E51').Value = '  +7.30%  '
